# Update the timetable entries on both sections to reflect the
# re-shuffled course schedule (sem3_timetable.xlsx).

$wb = $excel.ActiveWorkbook

# ----- Section_A -----
$wsA = $wb.Worksheets.Item("Section_A")

# Row 2 (9:00-10:30)
$wsA.Range("B2").Value = "CS264"
$wsA.Range("D2").Value = "CS263"
$wsA.Range("E2").Value = "Free"

# Row 3 (10:30-12:00)
$wsA.Range("B3").Value = "Free"
$wsA.Range("D3").Value = "Free"
$wsA.Range("E3").Value = "CS261"
$wsA.Range("F3").Value = "CS263"

# Row 5 (14:00-15:30)
$wsA.Range("D5").Value = "MA261"
$wsA.Range("E5").Value = "Free"
$wsA.Range("F5").Value = "CS264"

# Row 6 (15:30-17:00)
$wsA.Range("D6").Value = "Free"
$wsA.Range("E6").Value = "CS263"

# Row 7 (17:00-18:30)
$wsA.Range("C7").Value = "CS261"
$wsA.Range("F7").Value = "CS261"

# ----- Section_B -----
$wsB = $wb.Worksheets.Item("Section_B")

# Row 2 (9:00-10:30)
$wsB.Range("B2").Value = "Free"
$wsB.Range("C2").Value = "CS264"
$wsB.Range("E2").Value = "CS261"

# Row 3 (10:30-12:00)
$wsB.Range("D3").Value = "CS263"
$wsB.Range("E3").Value = "CS264"

# Row 5 (14:00-15:30)
$wsB.Range("B5").Value = "Free"
$wsB.Range("C5").Value = "CS263"
$wsB.Range("F5").Value = "CS261"

# Row 6 (15:30-17:00)
$wsB.Range("D6").Value = "Free"
$wsB.Range("E6").Value = "MA261"

# Row 7 (17:00-18:30)
$wsB.Range("B7").Value = "MA261"
$wsB.Range("C7").Value = "CS261"
$wsB.Range("E7").Value = "CS263"
